# "Cadmium and Tungsten Ore" — add two new materials (Cadmium, Tungsten) to
# the Materials sheet. This pushes the pre-existing "Stainless Steel" /
# "Bronze" summary rows down from rows 15/16 to rows 16/17, and their
# formulas are re-pointed at the (now shifted) Wood/Tin and Steel/Silver rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# --- Row 14: Cadmium (new literal data row) ---------------------------
$ws.Range("A14").Value = "Cadmium"
$ws.Range("B14").Value = 8.65
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 4.5999999999999996
$ws.Range("E14").Value = 1.46
$ws.Range("F14").Value = 25
$ws.Range("G14").Value = 78
$ws.Range("H14").Value = 321
$ws.Range("I14").Value = 767

# --- Row 15: Tungsten (new literal data row) ---------------------------
$ws.Range("A15").Value = "Tungsten"
$ws.Range("B15").Value = 19.25
$ws.Range("C15").Value = 7.5
$ws.Range("D15").Value = 40.5
$ws.Range("D15").NumberFormat = "0.0"
$ws.Range("E15").Value = 1.82
$ws.Range("F15").Value = 500
$ws.Range("G15").Value = 1000
$ws.Range("H15").Value = 3422
$ws.Range("I15").Value = 5000

# --- Row 16: Stainless Steel, now derived from Wood (row 5) + Tin (row 8)
$ws.Range("A16").Value = "Stainless Steel"
$ws.Range("B16").Formula = "=B5+B8"
$ws.Range("C16").Formula = "=C5+C8"
$ws.Range("D16").Formula = "=D5+D8"
$ws.Range("E16").Formula = "=E5+E8"
$ws.Range("F16").Formula = "=F5+F8"
$ws.Range("G16").Formula = "=G5+G8"
$ws.Range("H16").Formula = "=AVERAGE(H5,H8)"
$ws.Range("I16").Formula = "=AVERAGE(I5,I8)"

# --- Row 17: Bronze, now derived from Steel (row 4) + Silver (row 9) ---
# Row 17 used to be a blank placeholder row (style 3 on B:I); pull in the
# "data row" number formatting (style 9/8, like row 15/16) before writing
# the formulas so it renders the same way the summary rows above do.
$ws.Range("B16:I16").Copy()
$ws.Range("B17:I17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A17").Value = "Bronze"
$ws.Range("B17").Formula = "=B4+B9"
$ws.Range("C17").Formula = "=C4+C9"
$ws.Range("D17").Formula = "=D4+D9"
$ws.Range("E17").Formula = "=E4+E9"
$ws.Range("F17").Formula = "=F4+F9"
$ws.Range("G17").Formula = "=G4+G9"
$ws.Range("H17").Formula = "=AVERAGE(H4,H9)"
$ws.Range("I17").Formula = "=AVERAGE(I4,I9)"

# --- Update the saved selections on each sheet to match the edit session ---
$ws.Activate()
$ws.Range("F17").Select()

$wsPart1 = $wb.Worksheets.Item("Part 1")
$wsPart1.Activate()
$wsPart1.Range("N14").Select()

$wsPart2 = $wb.Worksheets.Item("Part 2")
$wsPart2.Activate()
$wsPart2.Range("B14").Select()

# Leave "Materials" as the active/selected tab, matching the saved file.
$ws.Activate()
